$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H17").Value = 2399.1177
$ws.Range("J17").Value = 2540.3845
$ws.Range("L17").Value = 7621.1535
$ws.Range("N17").Value = -7957.1535
$ws.Range("H40").Value = 5069.2144
$ws.Range("I40").Value = 4409.4443
$ws.Range("J40").Value = 6256.8
$ws.Range("K40").Value = 4409.4443
$ws.Range("L40").Value = 6256.8
$ws.Range("M40").Value = -4234.4443
$ws.Range("N40").Value = -6606.8
$ws.Range("H74").Value = 17976
$ws.Range("I74").Value = 17698.7
$ws.Range("K74").Value = 17698.7
$ws.Range("M74").Value = -16762.7
$ws.Range("H77").Value = 17976
$ws.Range("I77").Value = 17698.7
$ws.Range("K77").Value = 88493.5
$ws.Range("M77").Value = -83813.5
$ws.Range("H103").Value = 8931446
$ws.Range("I103").Value = 4425.75
$ws.Range("K103").Value = 13277.25
$ws.Range("M103").Value = -12691.25
$ws.Range("H106").Value = 2797.8125
$ws.Range("I106").Value = 2728.077
$ws.Range("K106").Value = 2728.077
$ws.Range("M106").Value = -2097.077
$ws.Range("H116").Value = 873060.6
$ws.Range("J116").Value = 1544.6666
$ws.Range("L116").Value = 1544.6666
$ws.Range("N116").Value = -8428.6666
$ws.Range("H118").Value = 1057.875
$ws.Range("J118").Value = 899.75
$ws.Range("L118").Value = 2699.25
$ws.Range("N118").Value = -6013.25
$ws.Range("H137").Value = 12932.777
$ws.Range("I137").Value = 2055.7144
$ws.Range("K137").Value = 6167.1432
$ws.Range("M137").Value = -3617.1432
$ws.Range("H138").Value = 278455.88
$ws.Range("J138").Value = 427500.1
$ws.Range("L138").Value = 1282500.3
$ws.Range("N138").Value = -1292780.3

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 872.63635
$ws.Range("I2").Value = 852
$ws.Range("K2").Value = 852
$ws.Range("M2").Value = -739
$ws.Range("H32").Value = 3840.4177
$ws.Range("I32").Value = 2971.5405
$ws.Range("K32").Value = 2971.5405
$ws.Range("M32").Value = -2684.5405
$ws.Range("H61").Value = 8916.429
$ws.Range("I61").Value = 0
$ws.Range("K61").Value = 0
$ws.Range("M61").ClearContents()
$ws.Range("H102").Value = 5864.9443
$ws.Range("I102").Value = 5785.5625
$ws.Range("K102").Value = 5785.5625
$ws.Range("M102").Value = -4163.5625
$ws.Range("H116").Value = 872.63635
$ws.Range("I116").Value = 852
$ws.Range("K116").Value = 852
$ws.Range("M116").Value = 1442
$ws.Range("H122").Value = 4676.154
$ws.Range("I122").Value = 4404.4688
$ws.Range("K122").Value = 13213.4064
$ws.Range("M122").Value = -10763.4064
$ws.Range("H136").Value = 8916.429
$ws.Range("I136").Value = 0
$ws.Range("K136").Value = 0
$ws.Range("M136").ClearContents()

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 872.63635
$ws.Range("I3").Value = 852
$ws.Range("K3").Value = 852
$ws.Range("M3").Value = -738
$ws.Range("H20").Value = 33339122
$ws.Range("I20").Value = 41673010
$ws.Range("K20").Value = 41673010
$ws.Range("M20").Value = -41672763
$ws.Range("H29").Value = 8200
$ws.Range("I29").Value = 8200
$ws.Range("K29").Value = 8200
$ws.Range("M29").Value = -7911
$ws.Range("H94").Value = 66667348
$ws.Range("I94").Value = 66667348
$ws.Range("K94").Value = 66667348
$ws.Range("M94").Value = -66666897
$ws.Range("H105").Value = 9632066
$ws.Range("I105").Value = 557685.4399999999
$ws.Range("K105").Value = 557685.4399999999
$ws.Range("M105").Value = -555938.4399999999
$ws.Range("H134").Value = 4250.8887
$ws.Range("J134").Value = 5126.75
$ws.Range("L134").Value = 15380.25
$ws.Range("N134").Value = -20450.25

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H7").Value = 62504400
$ws.Range("I7").Value = 5944
$ws.Range("K7").Value = 5944
$ws.Range("M7").Value = -5831
$ws.Range("H31").Value = 4304.8545
$ws.Range("I31").Value = 4029.6843
$ws.Range("J31").Value = 4450.0835
$ws.Range("K31").Value = 4029.6843
$ws.Range("L31").Value = 4450.0835
$ws.Range("M31").Value = -3734.6843
$ws.Range("N31").Value = -5040.0835
$ws.Range("H34").Value = 4304.8545
$ws.Range("I34").Value = 4029.6843
$ws.Range("J34").Value = 4450.0835
$ws.Range("K34").Value = 4029.6843
$ws.Range("L34").Value = 4450.0835
$ws.Range("M34").Value = -3827.6843
$ws.Range("N34").Value = -4854.0835
$ws.Range("H50").Value = 55130.2
$ws.Range("J50").Value = 55130.2
$ws.Range("L50").Value = 55130.2
$ws.Range("N50").Value = -56380.2
$ws.Range("H51").Value = 61419.5
$ws.Range("J51").Value = 61419.5
$ws.Range("L51").Value = 61419.5
$ws.Range("N51").Value = -62891.5
$ws.Range("H59").Value = 91917
$ws.Range("J59").Value = 91917
$ws.Range("L59").Value = 91917
$ws.Range("N59").Value = -94207
$ws.Range("H60").Value = 12833
$ws.Range("J60").Value = 17417.166
$ws.Range("L60").Value = 17417.166
$ws.Range("N60").Value = -18439.166
$ws.Range("H61").Value = 61419.5
$ws.Range("J61").Value = 61419.5
$ws.Range("L61").Value = 61419.5
$ws.Range("N61").Value = -62115.5
$ws.Range("H68").Value = 88448.5
$ws.Range("J68").Value = 88448.5
$ws.Range("L68").Value = 88448.5
$ws.Range("N68").Value = -89946.5
$ws.Range("H71").Value = 88448.5
$ws.Range("J71").Value = 88448.5
$ws.Range("L71").Value = 265345.5
$ws.Range("N71").Value = -272833.5
$ws.Range("H74").Value = 84656.5
$ws.Range("J74").Value = 84656.5
$ws.Range("L74").Value = 84656.5
$ws.Range("N74").Value = -86404.5
$ws.Range("H77").Value = 84656.5
$ws.Range("J77").Value = 84656.5
$ws.Range("L77").Value = 253969.5
$ws.Range("N77").Value = -262705.5

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H2").Value = 1343
$ws.Range("I2").Value = 60.333332
$ws.Range("K2").Value = 361.999992
$ws.Range("M2").Value = -248.999992
$ws.Range("H3").Value = 5488.5
$ws.Range("I3").Value = 5107
$ws.Range("K3").Value = 15321
$ws.Range("M3").Value = -15209
$ws.Range("H56").Value = 7252.75
$ws.Range("I56").Value = 7252.75
$ws.Range("K56").Value = 7252.75
$ws.Range("M56").Value = -6722.75
$ws.Range("H82").Value = 10744.167
$ws.Range("J82").Value = 10744.167
$ws.Range("L82").Value = 32232.501
$ws.Range("N82").Value = -33044.501
$ws.Range("H85").Value = 10744.167
$ws.Range("J85").Value = 10744.167
$ws.Range("L85").Value = 32232.501
$ws.Range("N85").Value = -35040.501
$ws.Range("H96").Value = 0
$ws.Range("J96").Value = 0
$ws.Range("L96").ClearContents()
$ws.Range("N96").Value = 0
$ws.Range("H122").Value = 1073.0454
$ws.Range("J122").Value = 1059.2222
$ws.Range("L122").Value = 9532.9998
$ws.Range("N122").Value = -14432.9998
$ws.Range("H140").Value = 9821.483
$ws.Range("I140").Value = 4748.5
$ws.Range("J140").Value = 22222.111
$ws.Range("K140").Value = 14245.5
$ws.Range("L140").Value = 66666.333
$ws.Range("M140").Value = -9065.5
$ws.Range("N140").Value = -77026.333

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H10").Value = 0
$ws.Range("I10").Value = 0
$ws.Range("K10").Value = 0
$ws.Range("M10").ClearContents()
$ws.Range("H20").Value = 25000
$ws.Range("J20").Value = 25000
$ws.Range("L20").Value = 25000
$ws.Range("N20").Value = -25490
$ws.Range("H69").Value = 25000
$ws.Range("J69").Value = 25000
$ws.Range("L69").Value = 25000
$ws.Range("N69").Value = -26498
$ws.Range("H70").Value = 16738296
$ws.Range("I70").Value = 23909328
$ws.Range("K70").Value = 23909328
$ws.Range("M70").Value = -23909058
$ws.Range("H72").Value = 25000
$ws.Range("J72").Value = 25000
$ws.Range("L72").Value = 75000
$ws.Range("N72").Value = -82488
$ws.Range("H73").Value = 16738296
$ws.Range("I73").Value = 23909328
$ws.Range("K73").Value = 23909328
$ws.Range("M73").Value = -23908392
$ws.Range("H102").Value = 3856.2432
$ws.Range("I102").Value = 916.069
$ws.Range("K102").Value = 916.069
$ws.Range("M102").Value = 705.931
$ws.Range("H107").Value = 997
$ws.Range("I107").Value = 997
$ws.Range("K107").Value = 997
$ws.Range("M107").Value = 923
$ws.Range("H122").Value = 12829179
$ws.Range("I122").Value = 38462536
$ws.Range("J122").Value = 12500
$ws.Range("K122").Value = 115387608
$ws.Range("L122").Value = 37500
$ws.Range("M122").Value = -115385158
$ws.Range("N122").Value = -42400
$ws.Range("H126").Value = 7026.5557
$ws.Range("I126").Value = 2856.25
$ws.Range("K126").Value = 8568.75
$ws.Range("M126").Value = -6098.75

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 4781.3335
$ws.Range("I7").Value = 4033.55
$ws.Range("J7").Value = 6917.857
$ws.Range("K7").Value = 4033.55
$ws.Range("L7").Value = 6917.857
$ws.Range("M7").Value = -3921.55
$ws.Range("N7").Value = -7141.857
$ws.Range("H40").Value = 5493.788
$ws.Range("I40").Value = 5623.893
$ws.Range("K40").Value = 5623.893
$ws.Range("M40").Value = -5487.893
$ws.Range("H55").Value = 368.8095
$ws.Range("I55").Value = 273.5
$ws.Range("K55").Value = 273.5
$ws.Range("M55").Value = -100.5
$ws.Range("H122").Value = 3479.6
$ws.Range("J122").Value = 3599.5
$ws.Range("L122").Value = 10798.5
$ws.Range("N122").Value = -15698.5
$ws.Range("H125").Value = 70000
$ws.Range("J125").Value = 70000
$ws.Range("L125").Value = 70000
$ws.Range("N125").Value = -79840
$ws.Range("H126").Value = 4781.3335
$ws.Range("I126").Value = 4033.55
$ws.Range("J126").Value = 6917.857
$ws.Range("K126").Value = 12100.65
$ws.Range("L126").Value = 20753.571
$ws.Range("M126").Value = -9630.650000000001
$ws.Range("N126").Value = -25693.571
$ws.Range("H127").Value = 69999
$ws.Range("J127").Value = 69999
$ws.Range("L127").Value = 69999
$ws.Range("N127").Value = -79919
$ws.Range("H130").Value = 0
$ws.Range("J130").Value = 0
$ws.Range("L130").ClearContents()
$ws.Range("N130").Value = 0
$ws.Range("H132").Value = 8994.066000000001
$ws.Range("J132").Value = 11244.25
$ws.Range("L132").Value = 33732.75
$ws.Range("N132").Value = -38792.75
$ws.Range("H136").Value = 4188
$ws.Range("I136").Value = 4051.5417
$ws.Range("K136").Value = 12154.6251
$ws.Range("M136").Value = -9604.625100000001

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H81").Value = 3655.2083
$ws.Range("I81").Value = 2780.5715
$ws.Range("J81").Value = 4879.7
$ws.Range("K81").Value = 5561.143
$ws.Range("L81").Value = 9759.4
$ws.Range("M81").Value = -4500.143
$ws.Range("N81").Value = -11881.4
$ws.Range("H84").Value = 3655.2083
$ws.Range("I84").Value = 2780.5715
$ws.Range("J84").Value = 4879.7
$ws.Range("K84").Value = 27805.715
$ws.Range("L84").Value = 48797
$ws.Range("M84").Value = -22501.715
$ws.Range("N84").Value = -59405
$ws.Range("H122").Value = 41670276
$ws.Range("I122").Value = 4174.5
$ws.Range("K122").Value = 12523.5
$ws.Range("M122").Value = -10073.5
$ws.Range("H126").Value = 2325
$ws.Range("I126").Value = 1322.125
$ws.Range("K126").Value = 3966.375
$ws.Range("M126").Value = -1496.375
$ws.Range("H132").Value = 10104248
$ws.Range("I132").Value = 11908070
$ws.Range("J132").Value = 2850
$ws.Range("K132").Value = 35724210
$ws.Range("L132").Value = 8550
$ws.Range("M132").Value = -35721680
$ws.Range("N132").Value = -13610
